$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps text formatting so values such as
# "1.001" or "1.000" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '28.578.71'
$ws.Range("E2").Value = '  +1.61%  '

# Row 3
$ws.Range("D3").Value = '1.827.84'
$ws.Range("E3").Value = '  +1.83%  '

# Row 4
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").Value = '317.77'

# Row 6
$ws.Range("E6").Value = '  -0.03%  '

# Row 7
$ws.Range("D7").Value = '0.5417'
$ws.Range("E7").Value = '  +0.19%  '

# Row 8
$ws.Range("D8").Value = '0.4021'
$ws.Range("E8").Value = '  +6.27%  '

# Row 9
$ws.Range("D9").Value = '0.07673'
$ws.Range("E9").Value = '  +3.02%  '

# Row 10
$ws.Range("D10").Value = '1.121'
$ws.Range("E10").Value = '  +2.50%  '

# Row 11
$ws.Range("D11").Value = '41.86'
$ws.Range("E11").Value = '  +0.27%  '

# Row 12
$ws.Range("D12").Value = '21.16'
$ws.Range("E12").Value = '  +3.05%  '

# Row 13
$ws.Range("D13").Value = '6.329'
$ws.Range("E13").Value = '  +3.51%  '

# Row 14
$ws.Range("D14").Value = '7.641'
$ws.Range("E14").Value = '  +5.48%  '

# Row 15
$ws.Range("D15").Value = '1.000'
$ws.Range("E15").Value = '  -0.05%  '

# Row 16
$ws.Range("D16").Value = '1.826.55'
$ws.Range("E16").Value = '  +1.73%  '

# Row 17
$ws.Range("D17").Value = '0.00001091'
$ws.Range("E17").Value = '  +2.95%  '

# Row 18
$ws.Range("D18").Value = '90.01'
$ws.Range("E18").Value = '  +1.06%  '

# Row 19
$ws.Range("D19").Value = '0.06601'
$ws.Range("E19").Value = '  +1.84%  '

# Row 20
$ws.Range("E20").Value = '  +3.12%  '

# Row 21
$ws.Range("E21").Value = '  -0.01%  '

# Row 22
$ws.Range("D22").Value = '6.068'
$ws.Range("E22").Value = '  +2.77%  '

# Row 23
$ws.Range("D23").Value = '28.589.49'
$ws.Range("E23").Value = '  +1.60%  '

# Row 24
$ws.Range("D24").Value = '11.19'
$ws.Range("E24").Value = '  +0.21%  '

# Row 25
$ws.Range("D25").Value = '2.276'
$ws.Range("E25").Value = '  +8.98%  '

# Row 26
$ws.Range("D26").Value = '157.97'
$ws.Range("E26").Value = '  +1.84%  '

# Row 27
$ws.Range("D27").Value = '2.458'
$ws.Range("E27").Value = '  +7.59%  '

# Row 28
$ws.Range("D28").Value = '20.76'
$ws.Range("E28").Value = '  +2.37%  '

# Row 29
$ws.Range("D29").Value = '2.035.91'
$ws.Range("E29").Value = '  +1.89%  '

# Row 30
$ws.Range("E30").Value = '  +2.42%  '

# Row 31
$ws.Range("E31").Value = '  +0.90%  '

# Row 32
$ws.Range("D32").Value = '0.1111'
$ws.Range("E32").Value = '  +4.93%  '

# Row 33
$ws.Range("D33").Value = '5.687'

# Row 34
$ws.Range("D34").Value = '0.07525'
$ws.Range("E34").Value = '  +16.03%  '

# Row 35
$ws.Range("D35").Value = '3.646'
$ws.Range("E35").Value = '  -0.23%  '

# Row 36
$ws.Range("D36").Value = '0.2248'
$ws.Range("E36").Value = '  -0.37%  '

# Row 37
$ws.Range("D37").Value = '0.02360'
$ws.Range("E37").Value = '  +2.98%  '

# Row 38
$ws.Range("E38").Value = '  +3.97%  '

# Row 39
$ws.Range("D39").Value = '8.886'
$ws.Range("E39").Value = '  +5.11%  '

# Row 40
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '11.38'
$ws.Range("E40").Value = '  +2.59%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.6303'
$ws.Range("E41").Value = '  +1.99%  '

# Row 42
$ws.Range("D42").Value = '1.189'
$ws.Range("E42").Value = '  +1.31%  '

# Row 43
$ws.Range("E43").Value = '  -0.04%  '

# Row 44
$ws.Range("D44").Value = '1.402'
$ws.Range("E44").Value = '  -3.34%  '

# Row 45
$ws.Range("D45").Value = '13.46'
$ws.Range("E45").Value = '  +1.36%  '

# Row 46
$ws.Range("D46").Value = '0.5888'
$ws.Range("E46").Value = '  +1.86%  '

# Row 47
$ws.Range("D47").Value = '3.709'
$ws.Range("E47").Value = '  +0.99%  '

# Row 48
$ws.Range("D48").Value = '125.12'
$ws.Range("E48").Value = '  +1.06%  '

# Row 49
$ws.Range("E49").Value = '  +4.15%  '

# Row 50
$ws.Range("E50").Value = '  +0.65%  '

# Row 51
$ws.Range("D51").Value = '0.06908'
$ws.Range("E51").Value = '  +1.38%  '
